$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row 18: "dry weights F2" update (new static-chamber entry) ---
# Copy the date-format style from row 17's date cell (A17) down to A18 first
# so the new cell inherits the existing m/d/yyyy format, then set the value.
$ws.Range("A17").Copy($ws.Range("A18"))
$ws.Range("A18").Value = 44881

$ws.Range("B18").Value = "F2"
$ws.Range("C18").Value = "Loligo"
$ws.Range("D18").Value = "Loligo"
# Leading apostrophe reproduces the "quote prefix" text entry as typed in
# the source file (value starts with ~').
$ws.Range("E18").Value = "'~'58 mL small static chambers"
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = "Y"
$ws.Range("H18").Value = "Y"
$ws.Range("I18").Value = "N"
$ws.Range("J18").Value = "N"

# Widen column E (Chamber_volume) so the longer note text fits.
$ws.Columns("E").ColumnWidth = 24.833333333333332

# --- View state left by the author after the edit ---
$excel.ActiveWindow.Zoom = 70
$ws.Range("C22").Select()
